$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full taxon-stats table rewritten: two additional PCB congeners
# (CB187, CB194) are now reported for every group (Bivalvia,
# Crustacea, Polychaeta, Actinopterygii). Rows 19-73 shift down
# accordingly and the used range grows from A1:F73 to A1:F81.
$data = @()
$data += ,@("type", "grp", "PCB", "min", "median", "max")
$data += ,@("Prey", "Bivalvia", "CB28", 58.72, 77.31, 102.37)
$data += ,@("Prey", "Bivalvia", "CB31", 27.67, 38.51, 59.39)
$data += ,@("Prey", "Bivalvia", "CB44", 126.19, 167.61, 235.3)
$data += ,@("Prey", "Bivalvia", "CB49", 194.4, 263.1, 382.9)
$data += ,@("Prey", "Bivalvia", "CB52", 284.75, 372.3, 526.5)
$data += ,@("Prey", "Bivalvia", "CB101", 794.35, 998.5599999999999, 1424.62)
$data += ,@("Prey", "Bivalvia", "CB105", 135.03, 170.26, 264.51)
$data += ,@("Prey", "Bivalvia", "CB110", 758.12, 972.48, 1452.46)
$data += ,@("Prey", "Bivalvia", "CB118", 611.5, 799.22, 1054.67)
$data += ,@("Prey", "Bivalvia", "CB128", 177.3, 263.93, 358.25)
$data += ,@("Prey", "Bivalvia", "CB132", 234.34, 309.94, 446.96)
$data += ,@("Prey", "Bivalvia", "CB138", 1108.02, 1450.09, 1993.91)
$data += ,@("Prey", "Bivalvia", "CB149", 832.29, 1161.32, 1539.87)
$data += ,@("Prey", "Bivalvia", "CB153", 2053.71, 2799.35, 3427.45)
$data += ,@("Prey", "Bivalvia", "CB156", 69.7, 97.8, 121.8)
$data += ,@("Prey", "Bivalvia", "CB170", 205.57, 313.07, 394.59)
$data += ,@("Prey", "Bivalvia", "CB180", 390.56, 653.3099999999999, 822.09)
$data += ,@("Prey", "Bivalvia", "CB187", 425.13, 807.77, 1113.15)
$data += ,@("Prey", "Bivalvia", "CB194", 83.12, 122.24, 162.06)
$data += ,@("Prey", "Crustacea", "CB28", 21.94, 39.11, 58.22)
$data += ,@("Prey", "Crustacea", "CB31", 12.4, 18.67, 35.95)
$data += ,@("Prey", "Crustacea", "CB44", 38.02, 50.67, 85.84)
$data += ,@("Prey", "Crustacea", "CB49", 41.14, 89.25, 91.34)
$data += ,@("Prey", "Crustacea", "CB52", 134.89, 140.89, 245.68)
$data += ,@("Prey", "Crustacea", "CB101", 62.08, 66.36, 100.48)
$data += ,@("Prey", "Crustacea", "CB105", 21.37, 26.63, 84.2)
$data += ,@("Prey", "Crustacea", "CB110", 48.9, 54.65, 76.56)
$data += ,@("Prey", "Crustacea", "CB118", 201.35, 371.9, 397.09)
$data += ,@("Prey", "Crustacea", "CB128", 10.9, 12.5, 22.58)
$data += ,@("Prey", "Crustacea", "CB132", 15.37, 18.37, 25.63)
$data += ,@("Prey", "Crustacea", "CB138", 173.53, 387.34, 495.54)
$data += ,@("Prey", "Crustacea", "CB149", 71.87, 90.45, 113.08)
$data += ,@("Prey", "Crustacea", "CB153", 347.47, 422.88, 517.02)
$data += ,@("Prey", "Crustacea", "CB156", 73.02, 84, 113.44)
$data += ,@("Prey", "Crustacea", "CB170", 99.47, 124.38, 240.12)
$data += ,@("Prey", "Crustacea", "CB180", 120.74, 159.56, 617.8099999999999)
$data += ,@("Prey", "Crustacea", "CB187", 99.59999999999999, 151.06, 688.5700000000001)
$data += ,@("Prey", "Crustacea", "CB194", 29.49, 49.65, 82.78)
$data += ,@("Prey", "Polychaeta", "CB28", 51.46, 67.89, 89.65000000000001)
$data += ,@("Prey", "Polychaeta", "CB31", 27.57, 37.04, 43)
$data += ,@("Prey", "Polychaeta", "CB44", 130.74, 154.84, 193.78)
$data += ,@("Prey", "Polychaeta", "CB49", 183.58, 229.12, 285.12)
$data += ,@("Prey", "Polychaeta", "CB52", 277.83, 335.19, 435.68)
$data += ,@("Prey", "Polychaeta", "CB101", 747.1799999999999, 849.26, 1006.43)
$data += ,@("Prey", "Polychaeta", "CB105", 111.27, 132.95, 204.3)
$data += ,@("Prey", "Polychaeta", "CB110", 550.5599999999999, 837.8200000000001, 964.75)
$data += ,@("Prey", "Polychaeta", "CB118", 563.98, 633.65, 899.46)
$data += ,@("Prey", "Polychaeta", "CB128", 193.32, 225.02, 272.47)
$data += ,@("Prey", "Polychaeta", "CB132", 190.34, 287.26, 326.16)
$data += ,@("Prey", "Polychaeta", "CB138", 1024.34, 1251.52, 1529.03)
$data += ,@("Prey", "Polychaeta", "CB149", 796.24, 1005.21, 1326.88)
$data += ,@("Prey", "Polychaeta", "CB153", 1784.42, 2007.96, 2743.61)
$data += ,@("Prey", "Polychaeta", "CB156", 57.1, 67.03, 104.78)
$data += ,@("Prey", "Polychaeta", "CB170", 222.53, 286.09, 356.74)
$data += ,@("Prey", "Polychaeta", "CB180", 483.31, 610.1, 726.83)
$data += ,@("Prey", "Polychaeta", "CB187", 426.66, 511.32, 650.88)
$data += ,@("Prey", "Polychaeta", "CB194", 77.39, 119.95, 145.84)
$data += ,@("Sole", "Actinopterygii", "CB28", 29.08, 39.44, 57.27)
$data += ,@("Sole", "Actinopterygii", "CB31", 4.47, 8.68, 18.47)
$data += ,@("Sole", "Actinopterygii", "CB44", 100.46, 138.7, 201.85)
$data += ,@("Sole", "Actinopterygii", "CB49", 148.17, 232.11, 283.35)
$data += ,@("Sole", "Actinopterygii", "CB52", 276.48, 427.9, 547.14)
$data += ,@("Sole", "Actinopterygii", "CB101", 772.6, 968.99, 1190.78)
$data += ,@("Sole", "Actinopterygii", "CB105", 124.15, 167.32, 208.69)
$data += ,@("Sole", "Actinopterygii", "CB110", 677.55, 825.8099999999999, 1005.34)
$data += ,@("Sole", "Actinopterygii", "CB118", 730.41, 958.64, 1117.26)
$data += ,@("Sole", "Actinopterygii", "CB128", 202.86, 255.76, 287.61)
$data += ,@("Sole", "Actinopterygii", "CB132", 164.32, 200.75, 233.8)
$data += ,@("Sole", "Actinopterygii", "CB138", 1158.09, 1473.97, 1796.05)
$data += ,@("Sole", "Actinopterygii", "CB149", 741.5, 853.86, 947.85)
$data += ,@("Sole", "Actinopterygii", "CB153", 2134.13, 2667.06, 2926.74)
$data += ,@("Sole", "Actinopterygii", "CB156", 67.06, 86.23, 90.56)
$data += ,@("Sole", "Actinopterygii", "CB170", 232.89, 281.08, 348.8)
$data += ,@("Sole", "Actinopterygii", "CB180", 499.99, 681.9, 783.8)
$data += ,@("Sole", "Actinopterygii", "CB187", 507.98, 599.96, 656.14)
$data += ,@("Sole", "Actinopterygii", "CB194", 79.56999999999999, 90.93000000000001, 104.4)
$data += ,@("Prey", "Bivalvia", "sumPCB", 9307.91, 11556.29, 15129.04)
$data += ,@("Prey", "Crustacea", "sumPCB", 2148.1, 2583.92, 3438.02)
$data += ,@("Prey", "Polychaeta", "sumPCB", 8011.82, 10624.91, 11784.82)
$data += ,@("Sole", "Actinopterygii", "sumPCB", 8654.84, 11494.79, 12596.26)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
